$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column for rows 2-5
# from 2023-10-25 (45224) to 2023-11-03 (45233)
$newValue = 45233

$ws.Range("C2").Value = $newValue
$ws.Range("C3").Value = $newValue
$ws.Range("C4").Value = $newValue
$ws.Range("C5").Value = $newValue
